# Apply the "Pull information and set up config excel file" edit.
$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# --- MASTER WORKSHEET (Sheet2) ---
# Remove the old 3rd column (Column C Header Text / Junk C3 text / formula column)
$ws2.Columns.Item(3).Delete()

# Remove the old junk formula rows 3-5 (1/2/3, 2/3/4, SUM rows)
$ws2.Range("A3:A5").EntireRow.Delete()

# Update header/junk row values
$ws2.Range("B2").Value = "Device Name 1"

# Add the new device rows in column B
$ws2.Range("B3").Value = "CAM1"
$ws2.Range("B4").Value = "CAM2"
$ws2.Range("B5").Value = "CAM3"
$ws2.Range("B6").Value = "GSS1"
$ws2.Range("B7").Value = "KEY1"
$ws2.Range("B8").Value = "KEY2"
$ws2.Range("B9").Value = "INP1"
$ws2.Range("B10").Value = "INP2"
$ws2.Range("B11").Value = "CAM4"

# Update selection on the sheet
$ws2.Range("I13").Select()
